# Update Object Pooling PPT Content
# - Set object to be found by game object name (add alt text / descr to the picture)
# - Adding and deleting resource data (re-saved picture dimensions shift slightly)

$p = $ppt.ActivePresentation

# Locate the picture named "그림 58" (shape id 1034) that PowerPoint placed on
# slide 8 of the deck. Search defensively by Id+Name instead of hard-coding
# slide/shape indices, since several slides reuse shape id 1034 for unrelated
# shapes.
$targetSlide = $null
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $shape = $sl.Shapes.Item($shi)
        if ($shape.Id -eq 1034 -and $shape.Name -eq "그림 58") {
            $targetSlide = $sl
            $targetShape = $shape
            break
        }
    }
    if ($targetShape -ne $null) { break }
}

if ($targetShape -eq $null) {
    # Fallback to the known location (8th slide, 3rd shape) if the search
    # above somehow fails to match.
    $targetSlide = $p.Slides.Item(8)
    $targetShape = $targetSlide.Shapes.Item(3)
}

# Give the picture a description pointing at the source resource file
# (this is surfaced in OOXML as the <p:cNvPr descr="..."/> attribute), so
# the object can be found by name/description like other game objects.
$targetShape.AlternativeText = "C:/Users/Admin1/AppData/Roaming/PolarisOffice/ETemp/14712_13418472/fImage115662033281.png"

# Nudge the picture's stored extents very slightly to match the re-saved
# resource data (was 2361565 x 3677920 EMU, now 2362200 x 3678555 EMU).
$targetShape.Width = 186.0
$targetShape.Height = 289.650035
